$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = @(
    62.31041666666667,
    62.22069444444443,
    62.1309722222222,
    62.04124999999999,
    61.95152777777776,
    61.86180555555555,
    61.77208333333331,
    61.68236111111111,
    61.59263888888887,
    61.50291666666666,
    61.41319444444443,
    61.32347222222222,
    61.23374999999999,
    61.14402777777778,
    61.05430555555554,
    60.96458333333334,
    60.8748611111111,
    60.78513888888887,
    60.69541666666666,
    60.60569444444442,
    60.51597222222222,
    60.42624999999998,
    60.33652777777777,
    60.24680555555554,
    60.15708333333333,
    60.0673611111111,
    59.97763888888889,
    59.88791666666665,
    59.79819444444445,
    59.70847222222221,
    59.61874999999998,
    59.52902777777777,
    59.43930555555553,
    59.34958333333333,
    59.25986111111109,
    59.17013888888889,
    59.08041666666665,
    58.99069444444444,
    58.90097222222221,
    58.81125,
    58.72152777777777,
    58.63180555555556,
    58.54208333333332,
    58.45236111111112,
    58.36263888888888,
    58.27291666666665,
    58.18319444444444,
    58.0934722222222,
    58.00375,
    57.91402777777776,
    57.82430555555555,
    57.73458333333332,
    57.64486111111111,
    57.55513888888888,
    57.46541666666667,
    57.37569444444443,
    57.28597222222223,
    57.19624999999999,
    57.10652777777779,
    57.01680555555555,
    56.92708333333331,
    56.83736111111111,
    56.74763888888887,
    56.65791666666667,
    56.56819444444443,
    56.47847222222222,
    56.38874999999999,
    56.29902777777778,
    56.20930555555555,
    56.11958333333334,
    56.0298611111111,
    55.9401388888889,
    55.85041666666666,
    55.76069444444445,
    55.67097222222222,
    55.58124999999998,
    55.49152777777778,
    55.40180555555554,
    55.31208333333333,
    55.2223611111111,
    55.13263888888889,
    55.04291666666666,
    54.95319444444445,
    54.86347222222221,
    54.77375000000001,
    54.68402777777777,
    54.59430555555556,
    54.50458333333333,
    54.41486111111109,
    54.32513888888889,
    54.23541666666665,
    54.14569444444444,
    54.05597222222221,
    53.96625,
    53.87652777777777,
    53.78680555555556,
    53.69708333333332,
    53.60736111111112,
    53.51763888888888,
    53.42791666666668,
    53.33819444444444,
    53.24847222222223,
    53.15875,
    53.06902777777776,
    52.97930555555556,
    52.88958333333332,
    52.79986111111111,
    52.71013888888888,
    52.62041666666667,
    52.53069444444444,
    52.44097222222223,
    52.35124999999999,
    52.26152777777779,
    52.17180555555555,
    52.08208333333334,
    51.99236111111111,
    51.9026388888889,
    51.81291666666667,
    51.72319444444443,
    51.63347222222222,
    51.54374999999999,
    51.45402777777778,
    51.36430555555555,
    51.27458333333334,
    51.1848611111111,
    51.0951388888889,
    51.00541666666666,
    50.91569444444445,
    50.82597222222222,
    50.73625000000001,
    50.64652777777778,
    50.55680555555554,
    50.46708333333333,
    50.3773611111111,
    50.28763888888889,
    50.19791666666666,
    50.10819444444445,
    50.01847222222221,
    49.92875000000001,
    49.83902777777777,
    49.74930555555557,
    49.65958333333333,
    49.56986111111112,
    49.48013888888889,
    49.39041666666668,
    49.30069444444445,
    49.21097222222221,
    49.12125,
    49.03152777777777,
    48.94180555555556,
    48.85208333333333,
    48.76236111111112,
    48.67263888888888,
    48.58291666666668,
    48.49319444444444,
    48.40347222222223,
    48.31375,
    48.22402777777779,
    48.13430555555556,
    48.04458333333335,
    47.95486111111111,
    47.86513888888888,
    47.77541666666667,
    47.68569444444444,
    47.59597222222223,
    47.50624999999999,
    47.41652777777779,
    47.32680555555555,
    47.23708333333335,
    47.14736111111111,
    47.0576388888889,
    46.96791666666667,
    46.87819444444446,
    46.78847222222223,
    46.69874999999999,
    46.60902777777778,
    46.51930555555555,
    46.42958333333334,
    46.33986111111111,
    46.2501388888889,
    46.16041666666666,
    46.07069444444446,
    45.98097222222222,
    45.89125000000001,
    45.80152777777778,
    45.71180555555557,
    45.62208333333334,
    45.53236111111113,
    45.44263888888889,
    45.35291666666666,
    45.26319444444445,
    45.17347222222222,
    45.08375000000001,
    44.99402777777777,
    44.90430555555557,
    44.81458333333333,
    44.72486111111112,
    44.63513888888889,
    44.54541666666668,
    44.45569444444445,
    45,
    44.8950611111111,
    44.79024444444444,
    44.68554999999999,
    44.58097777777778,
    44.47652777777778,
    44.3722,
    44.26799444444444,
    44.16391111111111,
    44.05994999999999,
    43.95611111111111,
    43.85239444444445,
    43.7488,
    43.64532777777777,
    43.54197777777778,
    43.43875,
    43.33564444444445,
    43.23266111111111,
    43.1298,
    43.02706111111112,
    42.92444444444445,
    42.82195,
    42.71957777777778,
    42.61732777777777,
    42.51520000000001,
    42.41319444444445,
    42.31131111111112,
    42.20955,
    42.10791111111111,
    42.00639444444445,
    41.90499999999999,
    41.80372777777778,
    41.70257777777778,
    41.60155,
    41.50064444444445,
    41.39986111111111,
    41.2992,
    41.19866111111111,
    41.09824444444445,
    40.99795,
    40.89777777777778,
    40.79772777777778,
    40.6978,
    40.59799444444445,
    40.49831111111111,
    40.39875000000001,
    40.29931111111112,
    40.19999444444445,
    40.10080000000001,
    40.00172777777778,
    39.90277777777778,
    39.80395000000001,
    39.70524444444445,
    39.60666111111112,
    39.50820000000001,
    39.40986111111112,
    39.31164444444445,
    39.21355000000001,
    39.11557777777778,
    39.01772777777779,
    38.92000000000001,
    38.82239444444446,
    38.72491111111111,
    38.62755000000001,
    38.53031111111112,
    38.43319444444445,
    38.33620000000001,
    38.23932777777779,
    38.14257777777778,
    38.04595000000001,
    37.94944444444445,
    37.85306111111111,
    37.75680000000001,
    37.66066111111112,
    37.56464444444445,
    37.46875000000001,
    37.37297777777778,
    37.27732777777779,
    37.18180000000001,
    37.08639444444446,
    36.99111111111112,
    36.89595000000001,
    36.80091111111112,
    36.70599444444446,
    36.61120000000001,
    36.51652777777778,
    36.42197777777778,
    36.32755000000002,
    36.23324444444445,
    36.13906111111113,
    36.04500000000001,
    35.95106111111112,
    35.85724444444445,
    35.76355000000002,
    35.66997777777779,
    35.57652777777779,
    35.48320000000001,
    35.38999444444446,
    35.29691111111112,
    35.20395000000001,
    35.11111111111112,
    35.01839444444446,
    34.92580000000001,
    34.8333277777778,
    34.74097777777779,
    34.64875000000001,
    34.55664444444446,
    34.46466111111113,
    34.37280000000001,
    34.28106111111113,
    34.18944444444446,
    34.09795000000001,
    34.00657777777779,
    33.91532777777779,
    33.82420000000001,
    33.73319444444446,
    33.64231111111113,
    33.55155000000001,
    33.46091111111112,
    33.37039444444446,
    33.28000000000002,
    33.18972777777779,
    33.0995777777778,
    33.00955000000001,
    32.91964444444446,
    32.82986111111112,
    32.74020000000002,
    32.65066111111113,
    32.56124444444445,
    32.47195000000001,
    32.3827777777778,
    32.29372777777779,
    32.20480000000001,
    32.11599444444446,
    32.02731111111112,
    31.93875000000001,
    31.85031111111113,
    31.76199444444446,
    31.67380000000002,
    31.58572777777779,
    31.4977777777778,
    31.40995000000002,
    31.32224444444446,
    31.23466111111113,
    31.14720000000002,
    31.05986111111113,
    30.97264444444446,
    30.88555000000002,
    30.79857777777779,
    30.71172777777779,
    30.62500000000002,
    30.53839444444446,
    30.45191111111113,
    30.36555000000002,
    30.27931111111113,
    30.19319444444446,
    30.10720000000002,
    30.0213277777778,
    29.9355777777778,
    29.84995000000002,
    29.76444444444446,
    29.67906111111113,
    29.59380000000002,
    29.50866111111113,
    29.42364444444446,
    29.33875000000002,
    29.25397777777779,
    29.1693277777778,
    29.08480000000002,
    29.00039444444446,
    28.91611111111113,
    28.83195000000002,
    28.74791111111113,
    28.66399444444446,
    28.58020000000002,
    28.4965277777778,
    28.4129777777778,
    28.32955000000002,
    28.24624444444447,
    28.16306111111113,
    28.08000000000002,
    27.99706111111113,
    27.91424444444447,
    27.83155000000002,
    27.7489777777778,
    27.6665277777778,
    27.58420000000002,
    27.50199444444446,
    27.41991111111113,
    27.33795000000002,
    27.25611111111114,
    27.17439444444447,
    27.09280000000002,
    27.0113277777778,
    26.9299777777778,
    26.84875000000002,
    26.76764444444447,
    26.68666111111114,
    26.60580000000002,
    26.52506111111113,
    26.44444444444447,
    26.36395000000002,
    26.2835777777778,
    26.2033277777778,
    26.12320000000002,
    26.04319444444447,
    25.96331111111114,
    25.88355000000002,
    25.80391111111113,
    25.72439444444447,
    25.64500000000002,
    25.56572777777781,
    25.4865777777778,
    25.40755000000003,
    25.32864444444447,
    25.24986111111113,
    25.17120000000002,
    25.09266111111113,
    25.01424444444447,
    24.93595000000003,
    24.8577777777778,
    24.7797277777778,
    24.70180000000002,
    24.62399444444447,
    24.54631111111113,
    24.46875000000002,
    24.39131111111114,
    24.31399444444447,
    24.23680000000002,
    24.1597277777778,
    24.0827777777778,
    24.00595000000003,
    23.92924444444447,
    23.85266111111114,
    23.77620000000003,
    23.69986111111114,
    23.62364444444447,
    23.54755000000003,
    23.47157777777781,
    23.3957277777778,
    23.32000000000003,
    23.24439444444447,
    23.16891111111114,
    23.09355000000003,
    23.01831111111114,
    22.94319444444447,
    22.86820000000003,
    22.7933277777778,
    22.7185777777778,
    22.64395000000003,
    22.56944444444447
)
$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 4).Value = $v
    $row = $row + 1
}
Write-Host "Updated $($values.Count) cells in column D"
